$d = $word.ActiveDocument

function Replace-ParagraphByFind($findText, $xml) {
  $rng = $d.Content
  $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
  if (-not $ok) {
    throw "Find failed for: $findText"
  }
  $paraRng = $rng.Paragraphs(1).Range
  $paraRng.InsertXML($xml)
}

Replace-ParagraphByFind '個原因，我們需要將天然' '<w:p w:rsidR="007E279C" w:rsidRDefault="007E279C" w:rsidP="007E279C"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r w:rsidR="00280B82"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>氣，酸氣是需要被脫除的，</w:t></w:r><w:r w:rsidR="00275328"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>因為以下幾</w:t></w:r><w:r w:rsidR="00280B82"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>個原因，我們需要將天然</w:t></w:r></w:p>'
Replace-ParagraphByFind '第一、提高熱值以及管線輸送能力，第二、因為液態水' '<w:p w:rsidR="007E279C" w:rsidRDefault="007E279C" w:rsidP="007E279C"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">    氣</w:t></w:r><w:r w:rsidR="00280B82"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>脫水，</w:t></w:r><w:r w:rsidR="00275328"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>第一、提高熱值以及管線輸送能力，第二、因為液態水</w:t></w:r></w:p>'
Replace-ParagraphByFind '在管線中移動易造成設備腐蝕，第三、' '<w:p w:rsidR="007E279C" w:rsidRDefault="007E279C" w:rsidP="007E279C"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r w:rsidR="00275328"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>在管線中移動易造成設備腐蝕，第三、</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>液態水在</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>冰點結冰，高於</w:t></w:r></w:p>'
Replace-ParagraphByFind '冰點時與天然氣中的氣體成分易形成固化水合物，' '<w:p w:rsidR="00E44E09" w:rsidRDefault="007E279C" w:rsidP="007E279C"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> 冰點時與天然氣中的氣體成分易形成固化水合物，</w:t></w:r><w:r w:rsidR="0028124F"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>而</w:t></w:r><w:r w:rsidR="0028124F" w:rsidRPr="0028124F"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>三乙二醇</w:t></w:r><w:r w:rsidR="00E44E09"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">因 </w:t></w:r></w:p>'
Replace-ParagraphByFind '為吸水性高，取得容易，是' '<w:p w:rsidR="00E44E09" w:rsidRDefault="003D21F2" w:rsidP="007E279C"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r><w:r w:rsidR="00E44E09"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> 為吸水性高，取得容易，是</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>常使用的脫水溶劑，</w:t></w:r><w:r w:rsidR="00E44E09"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>天然氣是極為重</w:t></w:r></w:p>'
Replace-ParagraphByFind '設備操作安全並可以有效率的運輸' '<w:p w:rsidR="00405544" w:rsidRDefault="00ED13D3" w:rsidP="00E44E09"><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shape id="_x0000_s1032" type="#_x0000_t32" style="position:absolute;margin-left:-2.3pt;margin-top:445.75pt;width:499.65pt;height:.55pt;flip:y;z-index:251663360" o:connectortype="straight"/></w:pict></w:r><w:r w:rsidR="00E44E09"><w:t xml:space="preserve">           </w:t></w:r><w:r w:rsidR="00E44E09" w:rsidRPr="00E44E09"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>保</w:t></w:r><w:r w:rsidR="00544AF7"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>設備操作安全並可以有效率的運輸。</w:t></w:r></w:p>'

Replace-ParagraphByFind '(六)' '<w:p><w:pPr><w:pStyle w:val="Standard"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:color w:val="002060"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:color w:val="002060"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>地質調查知識</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>地質調查是</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>石油探</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>勘</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>很重要的前置作業，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>雖然在我們的既定印象中地質調查就是到室外觀測後，把結果記錄下來就可以，不過影片課程介紹了一個很完整的觀測模式，首先需要進行事前的準備(資料蒐集等等)，再來才是到戶外進行觀測，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>最後還</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>必須回到室內，進行工作結果的彙整以及圖表繪製，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>影片中提到許多地形專有名詞，以及生成環境的介紹，因為以前的高中地理課都有學習到，所以覺得十分有親切感，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>閱讀起來也能吸收得比較快，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>我覺得地質學十分有趣，可以從一些地表的小特徵，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>推估地質的形成的原因，</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>沉積岩</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>的外觀還能推估地形的年齡，</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>地形探</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>勘</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>人員要忍受著外頭的日</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>曬</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>雨淋，日復一日的觀察，想必是對自己的專業很有熱忱才能堅持下來。</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>'

# Update the cached PAGE field result in the footer (3 -> 4)
$footer = $d.Sections(1).Footers(1)
$fld = $footer.Range.Fields(1)
$fld.Result.Text = "4"

Write-Host "Edit complete"
